# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2404"
#   "<name>_new" -> "<name>_FV2410"
# Then turn the used range into a native Excel Table ("Table1") with an
# AutoFilter, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Rename the header cells in row 1 (A1:U1): _old -> _FV2404, _new -> _FV2410
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    if ($null -eq $header) { continue }
    $header = [string]$header

    if ($header.EndsWith("_old")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2404"
    } elseif ($header.EndsWith("_new")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2410"
    }
}

# Convert the full data range (including the header row) into an Excel
# Table, which also adds the AutoFilter dropdowns on the header row.
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate()
$ws.Cells.Item(2, 1).Select()
$excel.ActiveWindow.FreezePanes = $true
